# Update "想去人数" (F column) values on sheet "展览" (exhibitions)
# and sheet "全部类型" (all types), per the authoritative diff.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value for "展览"
$exhibitionUpdates = @{
    4  = 53
    7  = 158
    9  = 28
    10 = 251
    15 = 856
    18 = 422
    20 = 67
    23 = 1279
    24 = 2931
    28 = 66
    29 = 1642
    32 = 25
    36 = 616
    38 = 12
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F-column value for "全部类型" (note: F15 differs from the
# "展览" sheet and is reset to 0 instead of 856).
$allTypesUpdates = @{
    4  = 53
    7  = 158
    9  = 28
    10 = 251
    15 = 0
    18 = 422
    20 = 67
    23 = 1279
    24 = 2931
    28 = 66
    29 = 1642
    32 = 25
    36 = 616
    38 = 12
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
